$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper to set a text-typed cell value without Excel re-interpreting
# numeric-looking strings (e.g. "1.00", "5.30") as numbers, which would
# silently drop the trailing zeros / formatting that the source data
# relies on. We flip the cell to Text just long enough to assign the
# value, then restore the "Normal" style so no stray number-format
# style is left behind on cells that didn't have one originally.
function Set-TextCell($row, $col, $value) {
    $cell = $ws.Cells.Item($row, $col)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

# row, D (Price), E (Volume(1h))
Set-TextCell 2  4 "68.575.57"
Set-TextCell 2  5 "  +0.52%  "

Set-TextCell 3  4 "2.704.08"
Set-TextCell 3  5 "  +2.26%  "

Set-TextCell 4  5 "  +0.07%  "

Set-TextCell 5  4 "599.37"
Set-TextCell 5  5 "  +0.48%  "

Set-TextCell 6  4 "159.87"
Set-TextCell 6  5 "  +2.24%  "

Set-TextCell 7  5 "  +0.04%  "

Set-TextCell 8  5 "  +0.01%  "

Set-TextCell 9  4 "2.704.68"
Set-TextCell 9  5 "  +2.32%  "

Set-TextCell 10 5 "  -3.81%  "

Set-TextCell 11 5 "  -0.76%  "

Set-TextCell 12 4 "5.30"
Set-TextCell 12 5 "  +0.96%  "

Set-TextCell 13 4 "0.360"
Set-TextCell 13 5 "  +2.46%  "

Set-TextCell 14 5 "  +1.04%  "

Set-TextCell 15 4 "3.191.73"
Set-TextCell 15 5 "  +2.05%  "

Set-TextCell 16 5 "  -1.99%  "

Set-TextCell 17 4 "68.570.18"
Set-TextCell 17 5 "  +0.62%  "

Set-TextCell 18 4 "2.688.65"
Set-TextCell 18 5 "  +1.24%  "

Set-TextCell 19 4 "11.91"
Set-TextCell 19 5 "  +4.76%  "

Set-TextCell 20 4 "366.85"
Set-TextCell 20 5 "  +1.20%  "

Set-TextCell 21 4 "7.64"
Set-TextCell 21 5 "  +2.66%  "

Set-TextCell 22 4 "4.57"
Set-TextCell 22 5 "  +3.69%  "

Set-TextCell 23 4 "4.92"
Set-TextCell 23 5 "  +1.98%  "

Set-TextCell 24 5 "  +3.34%  "

Set-TextCell 25 5 "  -0.41%  "

Set-TextCell 26 5 "  +0.01%  "

Set-TextCell 27 5 "  +4.16%  "

Set-TextCell 29 5 "  -0.54%  "

Set-TextCell 30 5 "  +0.33%  "

Set-TextCell 31 4 "578.21"
Set-TextCell 31 5 "  +3.62%  "

Set-TextCell 32 4 "8.29"
Set-TextCell 32 5 "  +3.48%  "

Set-TextCell 33 4 "1.44"
Set-TextCell 33 5 "  +2.42%  "

Set-TextCell 34 5 "  +4.95%  "

Set-TextCell 35 5 "  +5.68%  "

Set-TextCell 36 5 "  +2.11%  "

Set-TextCell 37 4 "1.00"
Set-TextCell 37 5 "  +0.04%  "

Set-TextCell 38 4 "20.10"
Set-TextCell 38 5 "  +3.44%  "

Set-TextCell 39 4 "161.50"
Set-TextCell 39 5 "  +0.25%  "

Set-TextCell 40 4 "0.382"
Set-TextCell 40 5 "  +2.38%  "

Set-TextCell 41 5 "  +1.96%  "

Set-TextCell 42 4 "5.44"
Set-TextCell 42 5 "  +1.99%  "

# Rows 43 and 44 swap places: WhiteBITCoin <-> dogwifhat
Set-TextCell 43 2 "dogwifhat"
Set-TextCell 43 3 "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
Set-TextCell 43 4 "2.66"
Set-TextCell 43 5 "  +1.39%  "

Set-TextCell 44 2 "WhiteBITCoin"
Set-TextCell 44 3 "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
Set-TextCell 44 4 "17.87"
Set-TextCell 44 5 "  +0.47%  "

Set-TextCell 45 5 "  +0.05%  "

Set-TextCell 46 4 "0.0₆0316"
Set-TextCell 46 5 "  -7.08%  "

Set-TextCell 47 4 "158.91"
Set-TextCell 47 5 "  +0.15%  "

Set-TextCell 48 5 "  +5.03%  "

Set-TextCell 49 4 "1.78"
Set-TextCell 49 5 "  +5.03%  "

Set-TextCell 50 5 "  +7.43%  "

Set-TextCell 51 4 "22.14"
Set-TextCell 51 5 "  +0.78%  "
